$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 80, shifting existing rows 80-206 down to 81-207.
$ws.Range("A80:R80").EntireRow.Insert()

# Populate the newly inserted row 80 with the new data point.
$ws.Range("A80").Value = 3
$ws.Range("B80").Value = "Femacal de La Calera"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44477
$ws.Range("E80").Value = 5
$ws.Range("F80").Value = 100112009
$ws.Range("G80").Value = "Acelga"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 320
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = 2094
$ws.Range("N80").Value = '$/docena de atados (6 kilos)'
$ws.Range("O80").Value = "Provincia de Quillota"
$ws.Range("P80").Value = 349
$ws.Range("Q80").Value = 6
$ws.Range("R80").Value = "Hortaliza"
